# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet
#    that surfaces it (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2) Narrow the per-language status columns (Overview E:F, zh-cn C,
#    de-de C) from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1) Update the status values -------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# --- 2) Narrow the status columns -------------------------------------
# Target stored (XML) column width is ~13.4102 characters. The Excel
# COM ColumnWidth property is expressed in "Normal"-style characters and
# is internally rounded to whole pixels before being persisted (stored
# width = pixel-snapped ColumnWidth + 5/6), so we request the
# ColumnWidth value whose pixel-snapped result lands on that target.
$targetColumnWidth = 12.576851254417766

$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = $targetColumnWidth
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $targetColumnWidth
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $targetColumnWidth
